# "Fixed Schematic and Updated PCB"
# Update the List of Materials worksheet: several Value/Particulars/Remarks
# cells in the Inductor & Capacitor rows are corrected, row 16 is made
# taller to fit the longer remark, and the active selection / scroll
# position of the sheet window is moved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (C5 ceramic cap): voltage rating corrected 1KV -> 400V
$ws.Range("E5").Value = "0.01uF 400V "

# Row 8 (C6 ceramic cap): Value column now specified
$ws.Range("E8").Value = "220nF 50V"

# Row 15 (L2 choke coil): Value/Particulars swapped with a Remarks note
$ws.Range("F15").Value = "7mH NTL-"
$ws.Range("G15").Value = "Choke Coil "

# Row 16 (L1 transformer-wound inductor): Value/Particulars corrected and
# the Remarks note clarifies it is a transformer-wound inductor
$ws.Range("E16").Value = "2.1mH"
$ws.Range("F16").Value = "2.1mH 200V"
$ws.Range("G16").Value = "Transformer wound type Inductor. Maintains power on falling edge of current signal - Flicker free continuous operation"
$ws.Rows.Item(16).RowHeight = 45.75

# Row 17 (MOV1 varistor): Particulars column now specified
$ws.Range("E17").Value = "320V/100pF"

# Move the window scroll position / active selection like the saved file
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("J12").Select()
